$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize the last name value in B3 (was "rai", now "Rai")
$ws.Range("B3").Value = "Rai"

# Update selection to D4 and scroll view back to top
$ws.Range("D4").Select()
